# "Latest data" refresh: the first 35 detail rows (rows 34-68, the
# "3566/2" .. "1865" particelle block) are no longer part of the
# particelle-non-trovate export, so the remaining data (previously rows
# 69-257) moves up to start at row 34, and the sheet shrinks from
# A1:C257 down to A1:C222.
#
# Column A is just a running index (row number - 2) and is left alone;
# only the codice_particella (B) / codice_comune_catastale (C) values
# need to shift up by 35 rows, after which the now-duplicated trailing
# 35 rows are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data block up: the values that used to live in B69:C257
# become the new B34:C222.
$ws.Range("B69:C257").Copy() | Out-Null
$ws.Range("B34").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

# Remove the now-redundant trailing rows so the sheet ends at row 222.
$ws.Rows("223:257").Delete() | Out-Null
